$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 308 (shifts existing rows 308..408 down to 309..409)
$ws.Rows.Item(308).Insert()

# Populate the newly inserted row 308 with the new weekly price record
$ws.Range("A308").Value = 10
$ws.Range("B308").Value = "Vega Modelo de Temuco"
$ws.Range("C308").Value = "La Araucanía"
$ws.Range("D308").Value = "2022-11-11"
$ws.Range("E308").Value = 9
$ws.Range("F308").Value = 100112009
$ws.Range("G308").Value = "Acelga"
$ws.Range("H308").Value = "Sin especificar"
$ws.Range("I308").Value = "Primera"
$ws.Range("J308").Value = 65
$ws.Range("K308").Value = 8000
$ws.Range("L308").Value = 8000
$ws.Range("M308").Value = 8000
$ws.Range("N308").Value = "$/docena de atados (12 kilos)"
$ws.Range("O308").Value = "Provincia de Cautín"
$ws.Range("P308").Value = 667
$ws.Range("Q308").Value = 12
$ws.Range("R308").Value = "Hortaliza"
